$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A4").Value = "firstName"
$ws.Range("B4").Value = "lastName"
$ws.Range("C4").Value = "dob"
$ws.Range("D4").Value = "age"
$ws.Range("E4").Value = "major"

$labels = $ws.Range("B4:E4")
$labels.NumberFormat = "General"
$labels.Font.Name = "Courier New"
$labels.Font.Size = 9
$labels.Font.Color = 11171480
$labels.HorizontalAlignment = 1
$labels.VerticalAlignment = -4108

$ws.Range("E4").Select()
